$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 411.1875
$ws.Cells.Item(2, 9).Value = 373.33334
$ws.Cells.Item(2, 10).Value = 524.75
$ws.Cells.Item(2, 11).Value = 373.33334
$ws.Cells.Item(2, 12).Value = 524.75
$ws.Cells.Item(2, 13).Value = -260.33334
$ws.Cells.Item(2, 14).Value = -750.75
$ws.Cells.Item(37, 8).Value = 78.57143000000001
$ws.Cells.Item(37, 9).Value = 78.57143000000001
$ws.Cells.Item(37, 11).Value = 235.71429
$ws.Cells.Item(37, 13).Value = -109.71429
$ws.Cells.Item(70, 8).Value = 4562.75
$ws.Cells.Item(70, 9).Value = 1700.5
$ws.Cells.Item(70, 11).Value = 5101.5
$ws.Cells.Item(70, 13).Value = -4831.5
$ws.Cells.Item(73, 8).Value = 4562.75
$ws.Cells.Item(73, 9).Value = 1700.5
$ws.Cells.Item(73, 11).Value = 5101.5
$ws.Cells.Item(73, 13).Value = -4165.5
$ws.Cells.Item(107, 8).Value = 773.0714
$ws.Cells.Item(107, 9).Value = 630.0909
$ws.Cells.Item(107, 10).Value = 1297.3334
$ws.Cells.Item(107, 11).Value = 630.0909
$ws.Cells.Item(107, 12).Value = 1297.3334
$ws.Cells.Item(107, 13).Value = 1289.9091
$ws.Cells.Item(107, 14).Value = -5137.3334
$ws.Cells.Item(132, 8).Value = 1937.1471
$ws.Cells.Item(132, 9).Value = 1619.4138
$ws.Cells.Item(132, 10).Value = 3780
$ws.Cells.Item(132, 11).Value = 4858.2414
$ws.Cells.Item(132, 12).Value = 11340
$ws.Cells.Item(132, 13).Value = -2328.2414
$ws.Cells.Item(132, 14).Value = -16400
$ws.Cells.Item(137, 8).Value = 1409.7894
$ws.Cells.Item(137, 9).Value = 1430.4375
$ws.Cells.Item(137, 11).Value = 4291.3125
$ws.Cells.Item(137, 13).Value = -1741.3125
$ws.Cells.Item(138, 8).Value = 3158.5823
$ws.Cells.Item(138, 10).Value = 4134.6
$ws.Cells.Item(138, 12).Value = 12403.8
$ws.Cells.Item(138, 14).Value = -22683.8

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3888.5273
$ws.Cells.Item(32, 9).Value = 3540
$ws.Cells.Item(32, 11).Value = 3540
$ws.Cells.Item(32, 13).Value = -3253
$ws.Cells.Item(61, 8).Value = 2876.5117
$ws.Cells.Item(61, 10).Value = 6225
$ws.Cells.Item(61, 12).Value = 6225
$ws.Cells.Item(61, 14).Value = -6649
$ws.Cells.Item(74, 8).Value = 6114
$ws.Cells.Item(74, 9).Value = 1362.4839
$ws.Cells.Item(74, 11).Value = 1362.4839
$ws.Cells.Item(74, 13).Value = -488.4838999999999
$ws.Cells.Item(77, 8).Value = 6114
$ws.Cells.Item(77, 9).Value = 1362.4839
$ws.Cells.Item(77, 11).Value = 6812.4195
$ws.Cells.Item(77, 13).Value = -2444.4195
$ws.Cells.Item(122, 8).Value = 2321.138
$ws.Cells.Item(122, 9).Value = 2152.8635
$ws.Cells.Item(122, 11).Value = 6458.5905
$ws.Cells.Item(122, 13).Value = -4008.5905
$ws.Cells.Item(136, 8).Value = 2876.5117
$ws.Cells.Item(136, 10).Value = 6225
$ws.Cells.Item(136, 12).Value = 18675
$ws.Cells.Item(136, 14).Value = -23775

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(140, 8).Value = 97000
$ws.Cells.Item(140, 10).Value = 97000
$ws.Cells.Item(140, 12).Value = 97000
$ws.Cells.Item(140, 14).Value = -107360

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4318.033
$ws.Cells.Item(132, 9).Value = 4197.9565
$ws.Cells.Item(132, 11).Value = 12593.8695
$ws.Cells.Item(132, 13).Value = -10063.8695

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 136.77777
$ws.Cells.Item(10, 9).Value = 147.28572
$ws.Cells.Item(10, 11).Value = 441.85716
$ws.Cells.Item(10, 13).Value = -302.85716
$ws.Cells.Item(107, 8).Value = 752.7222
$ws.Cells.Item(107, 9).Value = 327.4
$ws.Cells.Item(107, 10).Value = 916.3077
$ws.Cells.Item(107, 11).Value = 982.1999999999999
$ws.Cells.Item(107, 12).Value = 2748.9231
$ws.Cells.Item(107, 13).Value = 937.8000000000001
$ws.Cells.Item(107, 14).Value = -6588.9231
$ws.Cells.Item(108, 8).Value = 494.16666
$ws.Cells.Item(108, 9).Value = 494.16666
$ws.Cells.Item(108, 11).Value = 1482.49998
$ws.Cells.Item(108, 13).Value = 1397.50002

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 50151
$ws.Cells.Item(7, 9).Value = 302
$ws.Cells.Item(7, 11).Value = 302
$ws.Cells.Item(7, 13).Value = -190
$ws.Cells.Item(8, 8).Value = 50151
$ws.Cells.Item(8, 9).Value = 302
$ws.Cells.Item(8, 11).Value = 302
$ws.Cells.Item(8, 13).Value = -163
$ws.Cells.Item(10, 8).Value = 227011.33
$ws.Cells.Item(10, 9).Value = 401240.6
$ws.Cells.Item(10, 11).Value = 401240.6
$ws.Cells.Item(10, 13).Value = -401071.6
$ws.Cells.Item(11, 8).Value = 6174552
$ws.Cells.Item(11, 9).Value = 9252261
$ws.Cells.Item(11, 11).Value = 9252261
$ws.Cells.Item(11, 13).Value = -9252122
$ws.Cells.Item(33, 8).Value = 11754.75
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 13).ClearContents()
$ws.Cells.Item(102, 8).Value = 35715384
$ws.Cells.Item(102, 9).Value = 787
$ws.Cells.Item(102, 11).Value = 787
$ws.Cells.Item(102, 13).Value = 835
$ws.Cells.Item(141, 8).Value = 61302.125
$ws.Cells.Item(141, 10).Value = 61302.125
$ws.Cells.Item(141, 12).Value = 61302.125
$ws.Cells.Item(141, 14).Value = -71662.125

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(3, 8).Value = 269.33334
$ws.Cells.Item(3, 9).Value = 269.33334
$ws.Cells.Item(3, 11).Value = 269.33334
$ws.Cells.Item(3, 13).Value = -157.33334
$ws.Cells.Item(7, 8).Value = 12636
$ws.Cells.Item(7, 9).Value = 20001
$ws.Cells.Item(7, 11).Value = 20001
$ws.Cells.Item(7, 13).Value = -19889
$ws.Cells.Item(15, 8).Value = 269.33334
$ws.Cells.Item(15, 9).Value = 269.33334
$ws.Cells.Item(15, 11).Value = 269.33334
$ws.Cells.Item(15, 13).Value = -99.33334000000002
$ws.Cells.Item(46, 8).Value = 2337.8
$ws.Cells.Item(46, 9).Value = 1663
$ws.Cells.Item(46, 11).Value = 1663
$ws.Cells.Item(46, 13).Value = -1475
$ws.Cells.Item(108, 8).Value = 40000
$ws.Cells.Item(108, 10).Value = 40000
$ws.Cells.Item(108, 12).Value = 40000
$ws.Cells.Item(108, 14).Value = -47680
$ws.Cells.Item(126, 8).Value = 12636
$ws.Cells.Item(126, 9).Value = 20001
$ws.Cells.Item(126, 11).Value = 60003
$ws.Cells.Item(126, 13).Value = -57533
$ws.Cells.Item(129, 8).Value = 47136
$ws.Cells.Item(129, 10).Value = 47136
$ws.Cells.Item(129, 12).Value = 47136
$ws.Cells.Item(129, 14).Value = -57136
$ws.Cells.Item(136, 8).Value = 7371.3
$ws.Cells.Item(136, 9).Value = 6529.7144
$ws.Cells.Item(136, 10).Value = 9335
$ws.Cells.Item(136, 11).Value = 19589.1432
$ws.Cells.Item(136, 12).Value = 28005
$ws.Cells.Item(136, 13).Value = -17039.1432
$ws.Cells.Item(136, 14).Value = -33105

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5650
$ws.Cells.Item(62, 9).Value = 2350
$ws.Cells.Item(62, 11).Value = 2350
$ws.Cells.Item(62, 13).Value = -1726
$ws.Cells.Item(65, 8).Value = 5650
$ws.Cells.Item(65, 9).Value = 2350
$ws.Cells.Item(65, 11).Value = 11750
$ws.Cells.Item(65, 13).Value = -8630
$ws.Cells.Item(96, 8).Value = 1986.6364
$ws.Cells.Item(96, 9).Value = 1663.75
$ws.Cells.Item(96, 10).Value = 2171.1428
$ws.Cells.Item(96, 11).Value = 1663.75
$ws.Cells.Item(96, 12).Value = 2171.1428
$ws.Cells.Item(96, 13).Value = -290.75
$ws.Cells.Item(96, 14).Value = -4917.1428
$ws.Cells.Item(107, 8).Value = 673
$ws.Cells.Item(107, 9).Value = 673
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2019
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -99
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 5306.643
$ws.Cells.Item(126, 9).Value = 4724.375
$ws.Cells.Item(126, 10).Value = 6083
$ws.Cells.Item(126, 11).Value = 14173.125
$ws.Cells.Item(126, 12).Value = 18249
$ws.Cells.Item(126, 13).Value = -11703.125
$ws.Cells.Item(126, 14).Value = -23189
$ws.Cells.Item(129, 8).Value = 35000
$ws.Cells.Item(129, 9).Value = 20000
$ws.Cells.Item(129, 10).Value = 40000
$ws.Cells.Item(129, 11).Value = 20000
$ws.Cells.Item(129, 12).Value = 40000
$ws.Cells.Item(129, 13).Value = -15000
$ws.Cells.Item(129, 14).Value = -50000
$ws.Cells.Item(132, 8).Value = 2154.2666
$ws.Cells.Item(132, 9).Value = 2093.8572
$ws.Cells.Item(132, 11).Value = 6281.571599999999
$ws.Cells.Item(132, 13).Value = -3751.571599999999
$ws.Cells.Item(140, 8).Value = 74072.25
$ws.Cells.Item(140, 10).Value = 74072.25
$ws.Cells.Item(140, 12).Value = 74072.25
$ws.Cells.Item(140, 14).Value = -84432.25
